$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 107
$ws1.Range("F4").Value = 418
$ws1.Range("F7").Value = 1150
$ws1.Range("F8").Value = 385
$ws1.Range("F9").Value = 196
$ws1.Range("F12").Value = 377
$ws1.Range("F13").Value = 398
$ws1.Range("F14").Value = 788
$ws1.Range("F15").Value = 179
$ws1.Range("F16").Value = 723
$ws1.Range("F18").Value = 79
$ws1.Range("F19").Value = 1010
$ws1.Range("F20").Value = 460
$ws1.Range("F22").Value = 83
$ws1.Range("F24").Value = 28
$ws1.Range("F26").Value = 470

# Sheet "演出" (Performance) - column F ("想去人数")
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 365
$ws2.Range("F5").Value = 41
$ws2.Range("F10").Value = 629

# Sheet "全部类型" (All types) - column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 107
$ws4.Range("F6").Value = 418
$ws4.Range("F9").Value = 1150
$ws4.Range("F10").Value = 385
$ws4.Range("F11").Value = 196
$ws4.Range("F14").Value = 365
$ws4.Range("F16").Value = 41
$ws4.Range("F17").Value = 377
$ws4.Range("F20").Value = 398
$ws4.Range("F21").Value = 788
$ws4.Range("F22").Value = 179
$ws4.Range("F23").Value = 723
$ws4.Range("F25").Value = 79
$ws4.Range("F26").Value = 1010
$ws4.Range("F27").Value = 460
$ws4.Range("F31").Value = 83
$ws4.Range("F33").Value = 629
$ws4.Range("F35").Value = 28
$ws4.Range("F38").Value = 470
